# Applies the following textual corrections to the document:
#   1. Update the "last modified" date/time stamp.
#   2. Fix typo: vigileant -> vigilant
#   3. Fix typo: reasonnable -> reasonable
#   4. Fix typo: operating sysems -> operating systems

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceOne = 1

# 1. Update the date/time stamp in the Date-styled paragraph.
$d.Content.Find.Execute(
    "May  27, 2021 (11:54:01 PM)", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "May  28, 2021 (01:53:57 AM)", $wdReplaceOne
)

# 2. Typo fix: vigileant -> vigilant
$d.Content.Find.Execute(
    "vigileant", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "vigilant", $wdReplaceOne
)

# 3. Typo fix: reasonnable -> reasonable
$d.Content.Find.Execute(
    "reasonnable", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "reasonable", $wdReplaceOne
)

# 4. Typo fix: operating sysems -> operating systems
$d.Content.Find.Execute(
    "operating sysems", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "operating systems", $wdReplaceOne
)
